# Apply edits described by the commit "Atualizacoes 16 de janeiro de 2024."

$wb = $excel.ActiveWorkbook

# 1. Rename the second sheet from "Include from " to "Include from Stroke Type"
$wsInclude = $wb.Worksheets.Item("Include from ")
$wsInclude.Name = "Include from Stroke Type"

# 2. Update the Metadata sheet values
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version: 1.0.1 -> 0.0.0
$wsMeta.Range("B3").Value = "0.0.0"

# Title: "ValueSet for stroke type based on ICHOM" -> "Stroke type based on ICHOM"
$wsMeta.Range("B5").Value = "Stroke type based on ICHOM"

# Experimental: (empty) -> "false"
# Use a leading apostrophe so the engine stores it as text "false" instead of
# coercing it to the boolean FALSE; then restore the original cell formatting
# (the apostrophe trick flips on a "quote prefix" style) by re-pasting the
# formats from the neighbouring cell that already carries the right style.
$wsMeta.Range("B7").Value = "'false"
$wsMeta.Range("A7").Copy()
$wsMeta.Range("B7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Date: 2023-11-21T19:08:35-03:00 -> 2024-01-11T13:00:00-03:00
$wsMeta.Range("B8").Value = "2024-01-11T13:00:00-03:00"

# Description: "This ValueSet aims to categorize..." -> "ValueSet that aims to categorize..."
$wsMeta.Range("B12").Value = "ValueSet that aims to categorize the stroke type according to the ICHOM standard."

# 3. Update the Include sheet's System URI value
$wsInclude.Range("B7").Value = "https://molic-avc.gabriellesantosleandro.com/CodeSystem/StrokeTypeCS"
